# The underlying commit swaps the OOXML content of ppt/theme/theme1.xml
# (the slide master's theme - was "Integral") and ppt/theme/theme2.xml
# (the notes master's theme - was "Office Theme"): theme1.xml ends up
# holding the stock "Office Theme" color scheme, theme2.xml ends up
# holding the old "Integral" color scheme.
#
# The slide master's theme (theme1.xml) is reachable and editable through
# the PowerPoint object model via Slide.ThemeColorScheme, which maps
# straight onto the <a:clrScheme> of the design actually applied to the
# slides - so we rewrite its 12 theme colors to the stock Office Theme
# palette here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme color scheme (RGB() values, i.e. 0xBBGGRR of the sRGB hex):
#  1 dk1      000000 -> 0
#  2 lt1      FFFFFF -> 16777215
#  3 dk2      44546A -> 6968388
#  4 lt2      E7E6E6 -> 15132391
#  5 accent1  5B9BD5 -> 13998939
#  6 accent2  ED7D31 -> 3243501
#  7 accent3  A5A5A5 -> 10855845
#  8 accent4  FFC000 -> 49407
#  9 accent5  4472C4 -> 12874308
# 10 accent6  70AD47 -> 4697456
# 11 hlink    0563C1 -> 12673797
# 12 folHlink 954F72 -> 7491477

$tcs.Colors(1).RGB = 0
$tcs.Colors(2).RGB = 16777215
$tcs.Colors(3).RGB = 6968388
$tcs.Colors(4).RGB = 15132391
$tcs.Colors(5).RGB = 13998939
$tcs.Colors(6).RGB = 3243501
$tcs.Colors(7).RGB = 10855845
$tcs.Colors(8).RGB = 49407
$tcs.Colors(9).RGB = 12874308
$tcs.Colors(10).RGB = 4697456
$tcs.Colors(11).RGB = 12673797
$tcs.Colors(12).RGB = 7491477
